$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
}

$ws.Range("D2").Value = "27.248.31"
$ws.Range("E2").Value = "  +1.77%  "
$ws.Range("D3").Value = "1.818.53"
$ws.Range("E3").Value = "  +1.20%  "
$ws.Range("E4").Value = "  +0.17%  "
Set-TextValue "D5" "313.67"
$ws.Range("E5").Value = "  +1.63%  "
Set-TextValue "D6" "1.002"
$ws.Range("E6").Value = "  +0.25%  "
Set-TextValue "D7" "0.4651"
$ws.Range("E7").Value = "  +5.86%  "
Set-TextValue "D8" "0.3764"
$ws.Range("E8").Value = "  +2.33%  "
Set-TextValue "D9" "0.07411"
$ws.Range("E9").Value = "  +0.88%  "
Set-TextValue "D10" "0.8707"
$ws.Range("E10").Value = "  +1.74%  "
Set-TextValue "D11" "20.65"
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("D12").Value = "1.822.54"
$ws.Range("E12").Value = "  -6.71%  "
Set-TextValue "D13" "6.673"
$ws.Range("E13").Value = "  +0.80%  "
Set-TextValue "D14" "5.398"
$ws.Range("E14").Value = "  +2.55%  "
Set-TextValue "D15" "0.07096"
$ws.Range("E15").Value = "  +0.40%  "
Set-TextValue "D16" "92.18"
$ws.Range("E16").Value = "  +0.27%  "
$ws.Range("E17").Value = "  +0.19%  "
Set-TextValue "D18" "0.000008759"
$ws.Range("E18").Value = "  +1.43%  "
Set-TextValue "D19" "1.000"
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("E20").Value = "  +1.26%  "
$ws.Range("D21").Value = "27.247.11"
$ws.Range("E21").Value = "  +1.64%  "
Set-TextValue "D22" "5.321"
$ws.Range("E22").Value = "  +3.42%  "
Set-TextValue "D23" "10.91"
$ws.Range("E23").Value = "  +0.94%  "
$ws.Range("D24").Value = "2.048.27"
$ws.Range("E24").Value = "  +1.58%  "
Set-TextValue "D25" "1.933"
$ws.Range("E25").Value = "  -2.04%  "
Set-TextValue "D26" "151.65"
$ws.Range("E26").Value = "  +0.16%  "
Set-TextValue "D27" "2.265"
$ws.Range("E27").Value = "  +2.97%  "
$ws.Range("E28").Value = "  +1.34%  "
Set-TextValue "D29" "5.290"
$ws.Range("E29").Value = "  +2.36%  "
Set-TextValue "D30" "117.33"
$ws.Range("E30").Value = "  +0.13%  "
Set-TextValue "D31" "0.08908"
$ws.Range("E31").Value = "  +1.33%  "
Set-TextValue "D32" "0.7792"
$ws.Range("E32").Value = "  +5.66%  "
$ws.Range("E33").Value = "  +2.51%  "
Set-TextValue "D34" "4.530"
$ws.Range("E34").Value = "  +2.32%  "
Set-TextValue "D35" "2.922"
$ws.Range("E35").Value = "  +0.89%  "
Set-TextValue "D36" "1.001"
$ws.Range("E36").Value = "  +0.27%  "
Set-TextValue "D37" "1.104"
$ws.Range("E37").Value = "  +1.88%  "
Set-TextValue "D38" "0.01965"
$ws.Range("E38").Value = "  +0.67%  "
$ws.Range("E39").Value = "  +1.62%  "
Set-TextValue "D40" "7.282"
$ws.Range("E40").Value = "  +4.21%  "
$ws.Range("B41").Value = "MXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue "D41" "2.909"
$ws.Range("E41").Value = "  +3.12%  "
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D42" "2.374"
$ws.Range("E42").Value = "  +20.14%  "
Set-TextValue "D43" "0.5303"
$ws.Range("E43").Value = "  +1.37%  "
Set-TextValue "D44" "0.1691"
$ws.Range("E44").Value = "  +0.63%  "
Set-TextValue "D45" "8.617"
$ws.Range("E45").Value = "  +2.03%  "
Set-TextValue "D46" "0.5050"
$ws.Range("E46").Value = "  +0.20%  "
Set-TextValue "D47" "10.49"
$ws.Range("E47").Value = "  +0.98%  "
Set-TextValue "D48" "105.47"
$ws.Range("E48").Value = "  +0.50%  "
$ws.Range("E49").Value = "  +0.99%  "
Set-TextValue "D50" "1.001"
$ws.Range("E50").Value = "  +0.32%  "
Set-TextValue "D51" "0.06324"
$ws.Range("E51").Value = "  +0.66%  "
